$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 62.74008966666666
$ws.Cells.Item(2, 8).Value = 188.220269
$ws.Cells.Item(2, 9).Value = 0.6543216524118473
$ws.Cells.Item(2, 10).Value = 0.6543216524118471
$ws.Cells.Item(2, 13).Value = 145.7007446666667
$ws.Cells.Item(2, 14).Value = 437.1022340000001
$ws.Cells.Item(2, 15).Value = 0.2865937750105843
$ws.Cells.Item(2, 16).Value = 0.2865937750105843
$ws.Cells.Item(2, 17).Value = 9141.277784886772
$ws.Cells.Item(2, 18).Value = 82271.50006398094
$ws.Cells.Item(2, 19).Value = 0.1875245124358747
$ws.Cells.Item(2, 20).Value = 0.1875245124358746

$ws.Cells.Item(3, 7).Value = 62.74008966666666
$ws.Cells.Item(3, 8).Value = 188.220269
$ws.Cells.Item(3, 9).Value = 0.6543216524118473
$ws.Cells.Item(3, 10).Value = 0.6543216524118471
$ws.Cells.Item(3, 15).Value = 0.3320294904365841
$ws.Cells.Item(3, 16).Value = 0.3320294904365841
$ws.Cells.Item(3, 17).Value = 10590.50848101334
$ws.Cells.Item(3, 18).Value = 95314.57632912006
$ws.Cells.Item(3, 19).Value = 0.2172540848319294
$ws.Cells.Item(3, 20).Value = 0.2172540848319293

$ws.Cells.Item(4, 7).Value = 62.74008966666666
$ws.Cells.Item(4, 8).Value = 188.220269
$ws.Cells.Item(4, 9).Value = 0.6543216524118473
$ws.Cells.Item(4, 10).Value = 0.6543216524118471
$ws.Cells.Item(4, 13).Value = 128.1261546666667
$ws.Cells.Item(4, 14).Value = 384.378464
$ws.Cells.Item(4, 15).Value = 0.2520245069956105
$ws.Cells.Item(4, 16).Value = 0.2520245069956105
$ws.Cells.Item(4, 17).Value = 8038.646432431868
$ws.Cells.Item(4, 18).Value = 72347.81789188681
$ws.Cells.Item(4, 19).Value = 0.164905091865649
$ws.Cells.Item(4, 20).Value = 0.164905091865649

$ws.Cells.Item(5, 7).Value = 62.74008966666666
$ws.Cells.Item(5, 8).Value = 188.220269
$ws.Cells.Item(5, 9).Value = 0.6543216524118473
$ws.Cells.Item(5, 10).Value = 0.6543216524118471
$ws.Cells.Item(5, 13).Value = 65.761079
$ws.Cells.Item(5, 14).Value = 197.283237
$ws.Cells.Item(5, 15).Value = 0.1293522275572212
$ws.Cells.Item(5, 16).Value = 0.1293522275572212
$ws.Cells.Item(5, 17).Value = 4125.855993036749
$ws.Cells.Item(5, 18).Value = 37132.70393733074
$ws.Cells.Item(5, 19).Value = 0.08463796327839428
$ws.Cells.Item(5, 20).Value = 0.08463796327839426

$ws.Cells.Item(6, 9).Value = 0.1782000513806195
$ws.Cells.Item(6, 10).Value = 0.1782000513806195
$ws.Cells.Item(6, 13).Value = 145.7007446666667
$ws.Cells.Item(6, 14).Value = 437.1022340000001
$ws.Cells.Item(6, 15).Value = 0.2865937750105843
$ws.Cells.Item(6, 16).Value = 0.2865937750105843
$ws.Cells.Item(6, 17).Value = 2489.564826331038
$ws.Cells.Item(6, 18).Value = 22406.08343697934
$ws.Cells.Item(6, 19).Value = 0.05107102543225183
$ws.Cells.Item(6, 20).Value = 0.05107102543225182

$ws.Cells.Item(7, 9).Value = 0.1782000513806195
$ws.Cells.Item(7, 10).Value = 0.1782000513806195
$ws.Cells.Item(7, 15).Value = 0.3320294904365841
$ws.Cells.Item(7, 16).Value = 0.3320294904365841
$ws.Cells.Item(7, 19).Value = 0.05916767225568022
$ws.Cells.Item(7, 20).Value = 0.0591676722556802

$ws.Cells.Item(8, 9).Value = 0.1782000513806195
$ws.Cells.Item(8, 10).Value = 0.1782000513806195
$ws.Cells.Item(8, 13).Value = 128.1261546666667
$ws.Cells.Item(8, 14).Value = 384.378464
$ws.Cells.Item(8, 15).Value = 0.2520245069956105
$ws.Cells.Item(8, 16).Value = 0.2520245069956105
$ws.Cells.Item(8, 17).Value = 2189.270677517405
$ws.Cells.Item(8, 18).Value = 19703.43609765664
$ws.Cells.Item(8, 19).Value = 0.04491078009579309
$ws.Cells.Item(8, 20).Value = 0.04491078009579309

$ws.Cells.Item(9, 9).Value = 0.1782000513806195
$ws.Cells.Item(9, 10).Value = 0.1782000513806195
$ws.Cells.Item(9, 13).Value = 65.761079
$ws.Cells.Item(9, 14).Value = 197.283237
$ws.Cells.Item(9, 15).Value = 0.1293522275572212
$ws.Cells.Item(9, 16).Value = 0.1293522275572212
$ws.Cells.Item(9, 17).Value = 1123.648815896763
$ws.Cells.Item(9, 18).Value = 10112.83934307087
$ws.Cells.Item(9, 19).Value = 0.02305057359689441
$ws.Cells.Item(9, 20).Value = 0.02305057359689441

$ws.Cells.Item(10, 7).Value = 2.950144666666667
$ws.Cells.Item(10, 8).Value = 8.850434
$ws.Cells.Item(10, 9).Value = 0.03076730593473967
$ws.Cells.Item(10, 10).Value = 0.03076730593473966
$ws.Cells.Item(10, 13).Value = 145.7007446666667
$ws.Cells.Item(10, 14).Value = 437.1022340000001
$ws.Cells.Item(10, 15).Value = 0.2865937750105843
$ws.Cells.Item(10, 16).Value = 0.2865937750105843
$ws.Cells.Item(10, 17).Value = 429.8382748077285
$ws.Cells.Item(10, 18).Value = 3868.544473269556
$ws.Cells.Item(10, 19).Value = 0.008817718354742593
$ws.Cells.Item(10, 20).Value = 0.008817718354742592

$ws.Cells.Item(11, 7).Value = 2.950144666666667
$ws.Cells.Item(11, 8).Value = 8.850434
$ws.Cells.Item(11, 9).Value = 0.03076730593473967
$ws.Cells.Item(11, 10).Value = 0.03076730593473966
$ws.Cells.Item(11, 15).Value = 0.3320294904365841
$ws.Cells.Item(11, 16).Value = 0.3320294904365841
$ws.Cells.Item(11, 17).Value = 497.9835425569859
$ws.Cells.Item(11, 18).Value = 4481.851883012872
$ws.Cells.Item(11, 19).Value = 0.0102156529116181
$ws.Cells.Item(11, 20).Value = 0.0102156529116181

$ws.Cells.Item(12, 7).Value = 2.950144666666667
$ws.Cells.Item(12, 8).Value = 8.850434
$ws.Cells.Item(12, 9).Value = 0.03076730593473967
$ws.Cells.Item(12, 10).Value = 0.03076730593473966
$ws.Cells.Item(12, 13).Value = 128.1261546666667
$ws.Cells.Item(12, 14).Value = 384.378464
$ws.Cells.Item(12, 15).Value = 0.2520245069956105
$ws.Cells.Item(12, 16).Value = 0.2520245069956105
$ws.Cells.Item(12, 17).Value = 377.9906918503751
$ws.Cells.Item(12, 18).Value = 3401.916226653376
$ws.Cells.Item(12, 19).Value = 0.007754115109785884
$ws.Cells.Item(12, 20).Value = 0.007754115109785882

$ws.Cells.Item(13, 7).Value = 2.950144666666667
$ws.Cells.Item(13, 8).Value = 8.850434
$ws.Cells.Item(13, 9).Value = 0.03076730593473967
$ws.Cells.Item(13, 10).Value = 0.03076730593473966
$ws.Cells.Item(13, 13).Value = 65.761079
$ws.Cells.Item(13, 14).Value = 197.283237
$ws.Cells.Item(13, 15).Value = 0.1293522275572212
$ws.Cells.Item(13, 16).Value = 0.1293522275572212
$ws.Cells.Item(13, 17).Value = 194.0046964860953
$ws.Cells.Item(13, 18).Value = 1746.042268374858
$ws.Cells.Item(13, 19).Value = 0.003979819558593088
$ws.Cells.Item(13, 20).Value = 0.003979819558593087

$ws.Cells.Item(14, 7).Value = 13.108629
$ws.Cells.Item(14, 8).Value = 39.325887
$ws.Cells.Item(14, 9).Value = 0.1367109902727936
$ws.Cells.Item(14, 10).Value = 0.1367109902727935
$ws.Cells.Item(14, 13).Value = 145.7007446666667
$ws.Cells.Item(14, 14).Value = 437.1022340000001
$ws.Cells.Item(14, 15).Value = 0.2865937750105843
$ws.Cells.Item(14, 16).Value = 0.2865937750105843
$ws.Cells.Item(14, 17).Value = 1909.937006859062
$ws.Cells.Item(14, 18).Value = 17189.43306173156
$ws.Cells.Item(14, 19).Value = 0.03918051878771518
$ws.Cells.Item(14, 20).Value = 0.03918051878771516

$ws.Cells.Item(15, 7).Value = 13.108629
$ws.Cells.Item(15, 8).Value = 39.325887
$ws.Cells.Item(15, 9).Value = 0.1367109902727936
$ws.Cells.Item(15, 10).Value = 0.1367109902727935
$ws.Cells.Item(15, 15).Value = 0.3320294904365841
$ws.Cells.Item(15, 16).Value = 0.3320294904365841
$ws.Cells.Item(15, 17).Value = 2212.732677567644
$ws.Cells.Item(15, 18).Value = 19914.5940981088
$ws.Cells.Item(15, 19).Value = 0.04539208043735646
$ws.Cells.Item(15, 20).Value = 0.04539208043735645

$ws.Cells.Item(16, 7).Value = 13.108629
$ws.Cells.Item(16, 8).Value = 39.325887
$ws.Cells.Item(16, 9).Value = 0.1367109902727936
$ws.Cells.Item(16, 10).Value = 0.1367109902727935
$ws.Cells.Item(16, 13).Value = 128.1261546666667
$ws.Cells.Item(16, 14).Value = 384.378464
$ws.Cells.Item(16, 15).Value = 0.2520245069956105
$ws.Cells.Item(16, 16).Value = 0.2520245069956105
$ws.Cells.Item(16, 17).Value = 1679.558226721952
$ws.Cells.Item(16, 18).Value = 15116.02404049757
$ws.Cells.Item(16, 19).Value = 0.0344545199243825
$ws.Cells.Item(16, 20).Value = 0.03445451992438249

$ws.Cells.Item(17, 7).Value = 13.108629
$ws.Cells.Item(17, 8).Value = 39.325887
$ws.Cells.Item(17, 9).Value = 0.1367109902727936
$ws.Cells.Item(17, 10).Value = 0.1367109902727935
$ws.Cells.Item(17, 13).Value = 65.761079
$ws.Cells.Item(17, 14).Value = 197.283237
$ws.Cells.Item(17, 15).Value = 0.1293522275572212
$ws.Cells.Item(17, 16).Value = 0.1293522275572212
$ws.Cells.Item(17, 17).Value = 862.037587250691
$ws.Cells.Item(17, 18).Value = 7758.338285256219
$ws.Cells.Item(17, 19).Value = 0.01768387112333945
$ws.Cells.Item(17, 20).Value = 0.01768387112333945
